$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 445 - this shifts existing rows 445..496 down to 446..497
# and copies formatting (e.g. the date-style on column D) from the row above, matching
# the target workbook where D445 keeps style s="2".
$ws.Rows.Item(445).Insert()

# Fill in the new row 445 with this week's record (same dimension/location/category
# metadata as the surrounding weekly entries, with the new week's figures).
$ws.Cells.Item(445, 1).Value = 3
$ws.Cells.Item(445, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(445, 3).Value = "Coquimbo"
$ws.Cells.Item(445, 4).Value = 44776
$ws.Cells.Item(445, 5).Value = 5
$ws.Cells.Item(445, 6).Value = 100112003
$ws.Cells.Item(445, 7).Value = "Ajo"
$ws.Cells.Item(445, 8).Value = "Chino"
$ws.Cells.Item(445, 9).Value = "Primera"
$ws.Cells.Item(445, 10).Value = 85
$ws.Cells.Item(445, 11).Value = 26000
$ws.Cells.Item(445, 12).Value = 27000
$ws.Cells.Item(445, 13).Value = 26471
$ws.Cells.Item(445, 14).Value = "$/caja 10 kilos"
$ws.Cells.Item(445, 15).Value = "China"
$ws.Cells.Item(445, 16).Value = 2647
$ws.Cells.Item(445, 17).Value = 10
$ws.Cells.Item(445, 18).Value = "Hortaliza"
